# Rename the worksheet ("Anatomy" -> "Attendance Log")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Attendance Log"

# Delete row 4 first (old last data row), which shifts nothing but shrinks the
# used range / dimension from A1:F4 down to A1:F3, matching the target sheet.
$ws.Rows.Item(4).Delete()

# Ensure the ID columns keep their text ("number stored as text") semantics
# before writing the new ID values, same way the original values were stored.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A3").NumberFormat = "@"

# Row 2: replace the old row-2 values with the old row-3 id and the new
# attendance-log field values
$ws.Range("A2").Value = "424346"
$ws.Range("B2").Value = "Unknown"
$ws.Range("C2").Value = "14/08/2025"
$ws.Range("D2").Value = "9:31:23 AM"
$ws.Range("E2").Value = "QR Scan"
$ws.Range("F2").Value = "Unknown"

# Row 3: replace the old row-3 values (now shifted, previously old row-4) with
# the old row-4 id and the new attendance-log field values
$ws.Range("A3").Value = "676767"
$ws.Range("B3").Value = "Unknown"
$ws.Range("C3").Value = "14/08/2025"
$ws.Range("D3").Value = "9:31:25 AM"
$ws.Range("E3").Value = "QR Scan"
$ws.Range("F3").Value = "Unknown"

# Set explicit column widths to match the target worksheet
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 15
$ws.Columns.Item(6).ColumnWidth = 25
